$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting D:K -> E:L
$ws.Columns("D:D").Insert()

# Populate the new column D with FY2018 data (and NA/blank placeholders)
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 31900
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = 90700
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 145700
$ws.Range("D18").Value = -113700
$ws.Range("D20").Value = 3800
$ws.Range("D21").Value = -106700
$ws.Range("D22").Value = "NA"
$ws.Range("D23").Value = -110000
$ws.Range("D24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -110000
$ws.Range("D27").Value = -110000
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -3800
$ws.Range("D33").Value = -110000
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -110000
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 134800
$ws.Range("D42").Value = 234200
$ws.Range("D43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 5800
$ws.Range("D46").Value = 374800
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 40200
$ws.Range("D49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 5400
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 420400
$ws.Range("D57").Value = 5300
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 30600
$ws.Range("D60").Value = 35900
$ws.Range("D61").Value = 32400
$ws.Range("D62").Value = 115900
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 184200
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -416300
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 236200
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = -110000
$ws.Range("D83").Value = 3300
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = -45700
$ws.Range("D91").Value = -4800
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -53100
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 86900
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -11900

# Apply number formatting + font to match the existing column style (s=3 data cells)
$ws.Range("D7:D35").NumberFormat = "#,##0"
$ws.Range("D7:D35").Font.Name = "Verdana"
$ws.Range("D7:D35").Font.Size = 12
$ws.Range("D38:D77").NumberFormat = "#,##0"
$ws.Range("D38:D77").Font.Name = "Verdana"
$ws.Range("D38:D77").Font.Size = 12
$ws.Range("D80:D102").NumberFormat = "#,##0"
$ws.Range("D80:D102").Font.Name = "Verdana"
$ws.Range("D80:D102").Font.Size = 12

# Period-ending header rows (7, 38, 80) use the bold date style (s=2)
$ws.Range("D7").NumberFormat = "[$-409]d\-mmm\-yy;@"
$ws.Range("D7").Font.Bold = $true
$ws.Range("D38").NumberFormat = "[$-409]d\-mmm\-yy;@"
$ws.Range("D38").Font.Bold = $true
$ws.Range("D80").NumberFormat = "[$-409]d\-mmm\-yy;@"
$ws.Range("D80").Font.Bold = $true
